$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 223
$ws.Cells.Item(223, 4).Value = 44782
$ws.Cells.Item(223, 10).Value = 30

# Row 224
$ws.Cells.Item(224, 4).Value = 44200
$ws.Cells.Item(224, 11).Value = 4000
$ws.Cells.Item(224, 12).Value = 4000
$ws.Cells.Item(224, 13).Value = 4000
$ws.Cells.Item(224, 16).Value = 1333

# Row 225
$ws.Cells.Item(225, 4).Value = 44424
$ws.Cells.Item(225, 10).Value = 20
$ws.Cells.Item(225, 11).Value = 8000
$ws.Cells.Item(225, 12).Value = 8000
$ws.Cells.Item(225, 13).Value = 8000
$ws.Cells.Item(225, 16).Value = 2667

# Row 226
$ws.Cells.Item(226, 4).Value = 44259
$ws.Cells.Item(226, 10).Value = 80
$ws.Cells.Item(226, 11).Value = 5000
$ws.Cells.Item(226, 12).Value = 5000
$ws.Cells.Item(226, 13).Value = 5000
$ws.Cells.Item(226, 16).Value = 1667

# Row 227
$ws.Cells.Item(227, 4).Value = 44202
$ws.Cells.Item(227, 10).Value = 20
$ws.Cells.Item(227, 11).Value = 4000
$ws.Cells.Item(227, 12).Value = 4000
$ws.Cells.Item(227, 13).Value = 4000
$ws.Cells.Item(227, 16).Value = 1333

# Row 228
$ws.Cells.Item(228, 4).Value = 44371

# Row 229
$ws.Cells.Item(229, 4).Value = 44162
$ws.Cells.Item(229, 10).Value = 40

# Row 230
$ws.Cells.Item(230, 4).Value = 44235
$ws.Cells.Item(230, 10).Value = 110

# Row 231
$ws.Cells.Item(231, 4).Value = 44662
$ws.Cells.Item(231, 10).Value = 20
$ws.Cells.Item(231, 11).Value = 5000
$ws.Cells.Item(231, 12).Value = 5000
$ws.Cells.Item(231, 13).Value = 5000
$ws.Cells.Item(231, 16).Value = 1667

# Row 232
$ws.Cells.Item(232, 4).Value = 44427
$ws.Cells.Item(232, 10).Value = 30
$ws.Cells.Item(232, 11).Value = 8000
$ws.Cells.Item(232, 12).Value = 8000
$ws.Cells.Item(232, 13).Value = 8000
$ws.Cells.Item(232, 16).Value = 2667

# Row 233
$ws.Cells.Item(233, 4).Value = 44441
$ws.Cells.Item(233, 11).Value = 10000
$ws.Cells.Item(233, 12).Value = 10000
$ws.Cells.Item(233, 13).Value = 10000
$ws.Cells.Item(233, 16).Value = 3333

# Row 234
$ws.Cells.Item(234, 4).Value = 44708
$ws.Cells.Item(234, 10).Value = 20
$ws.Cells.Item(234, 11).Value = 6000
$ws.Cells.Item(234, 12).Value = 6000
$ws.Cells.Item(234, 13).Value = 6000
$ws.Cells.Item(234, 16).Value = 2000

# Row 235
$ws.Cells.Item(235, 4).Value = 44174
$ws.Cells.Item(235, 10).Value = 30
$ws.Cells.Item(235, 11).Value = 4000
$ws.Cells.Item(235, 12).Value = 4000
$ws.Cells.Item(235, 13).Value = 4000
$ws.Cells.Item(235, 16).Value = 1333

# Row 236
$ws.Cells.Item(236, 4).Value = 44419
$ws.Cells.Item(236, 10).Value = 65
$ws.Cells.Item(236, 11).Value = 10000
$ws.Cells.Item(236, 12).Value = 10000
$ws.Cells.Item(236, 13).Value = 10000
$ws.Cells.Item(236, 16).Value = 3333

# Row 237
$ws.Cells.Item(237, 4).Value = 44280
$ws.Cells.Item(237, 10).Value = 95
$ws.Cells.Item(237, 11).Value = 4000
$ws.Cells.Item(237, 12).Value = 4000
$ws.Cells.Item(237, 13).Value = 4000
$ws.Cells.Item(237, 16).Value = 1333

# Row 238
$ws.Cells.Item(238, 4).Value = 44518
$ws.Cells.Item(238, 10).Value = 40
$ws.Cells.Item(238, 15).Value = "Provincia de Cautín"

# Row 239
$ws.Cells.Item(239, 4).Value = 44412
$ws.Cells.Item(239, 11).Value = 5000
$ws.Cells.Item(239, 12).Value = 5000
$ws.Cells.Item(239, 13).Value = 5000
$ws.Cells.Item(239, 15).Value = "Región Metropolitana"
$ws.Cells.Item(239, 16).Value = 1667

# Row 240
$ws.Cells.Item(240, 4).Value = 44483
$ws.Cells.Item(240, 10).Value = 30
$ws.Cells.Item(240, 11).Value = 6000
$ws.Cells.Item(240, 12).Value = 6000
$ws.Cells.Item(240, 13).Value = 6000
$ws.Cells.Item(240, 16).Value = 2000

# Row 241
$ws.Cells.Item(241, 4).Value = 44175
$ws.Cells.Item(241, 10).Value = 40
$ws.Cells.Item(241, 11).Value = 4000
$ws.Cells.Item(241, 12).Value = 4000
$ws.Cells.Item(241, 13).Value = 4000
$ws.Cells.Item(241, 16).Value = 1333

# Row 242
$ws.Cells.Item(242, 4).Value = 44469
$ws.Cells.Item(242, 10).Value = 30
$ws.Cells.Item(242, 11).Value = 7000
$ws.Cells.Item(242, 12).Value = 7000
$ws.Cells.Item(242, 13).Value = 7000
$ws.Cells.Item(242, 16).Value = 2333

# Row 243
$ws.Cells.Item(243, 4).Value = 44434
$ws.Cells.Item(243, 10).Value = 45
$ws.Cells.Item(243, 11).Value = 8000
$ws.Cells.Item(243, 12).Value = 8000
$ws.Cells.Item(243, 13).Value = 8000
$ws.Cells.Item(243, 16).Value = 2667

# Row 244
$ws.Cells.Item(244, 4).Value = 44253
$ws.Cells.Item(244, 10).Value = 65
$ws.Cells.Item(244, 11).Value = 5000
$ws.Cells.Item(244, 12).Value = 5000
$ws.Cells.Item(244, 13).Value = 5000
$ws.Cells.Item(244, 16).Value = 1667

# Row 245
$ws.Cells.Item(245, 10).Value = 20
$ws.Cells.Item(245, 11).Value = 6000
$ws.Cells.Item(245, 12).Value = 6000
$ws.Cells.Item(245, 13).Value = 6000
$ws.Cells.Item(245, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(245, 16).Value = 2000

# Row 246
$ws.Cells.Item(246, 4).Value = 44494
$ws.Cells.Item(246, 10).Value = 30
$ws.Cells.Item(246, 11).Value = 2000
$ws.Cells.Item(246, 12).Value = 2000
$ws.Cells.Item(246, 13).Value = 2000
$ws.Cells.Item(246, 15).Value = "Región Metropolitana"
$ws.Cells.Item(246, 16).Value = 667

# Row 247
$ws.Cells.Item(247, 10).Value = 105
$ws.Cells.Item(247, 11).Value = 5000
$ws.Cells.Item(247, 12).Value = 6000
$ws.Cells.Item(247, 13).Value = 5619
$ws.Cells.Item(247, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(247, 16).Value = 1873

# Row 248
$ws.Cells.Item(248, 4).Value = 44487
$ws.Cells.Item(248, 10).Value = 75
$ws.Cells.Item(248, 11).Value = 2000
$ws.Cells.Item(248, 12).Value = 2000
$ws.Cells.Item(248, 13).Value = 2000
$ws.Cells.Item(248, 15).Value = "Región Metropolitana"
$ws.Cells.Item(248, 16).Value = 667

# Row 249
$ws.Cells.Item(249, 10).Value = 10
$ws.Cells.Item(249, 11).Value = 5000
$ws.Cells.Item(249, 12).Value = 5000
$ws.Cells.Item(249, 13).Value = 5000
$ws.Cells.Item(249, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(249, 16).Value = 1667

# Row 250
$ws.Cells.Item(250, 4).Value = 44356
$ws.Cells.Item(250, 10).Value = 20
$ws.Cells.Item(250, 11).Value = 2500
$ws.Cells.Item(250, 12).Value = 2500
$ws.Cells.Item(250, 13).Value = 2500
$ws.Cells.Item(250, 15).Value = "Región Metropolitana"
$ws.Cells.Item(250, 16).Value = 833

# Row 251
$ws.Cells.Item(251, 4).Value = 44301
$ws.Cells.Item(251, 10).Value = 80
$ws.Cells.Item(251, 11).Value = 7000
$ws.Cells.Item(251, 12).Value = 7000
$ws.Cells.Item(251, 13).Value = 7000
$ws.Cells.Item(251, 16).Value = 2333

# Row 252
$ws.Cells.Item(252, 4).Value = 44484
$ws.Cells.Item(252, 10).Value = 35
$ws.Cells.Item(252, 11).Value = 5000
$ws.Cells.Item(252, 13).Value = 5571
$ws.Cells.Item(252, 16).Value = 1857

# Row 253
$ws.Cells.Item(253, 10).Value = 40
$ws.Cells.Item(253, 11).Value = 6000
$ws.Cells.Item(253, 12).Value = 6000
$ws.Cells.Item(253, 13).Value = 6000
$ws.Cells.Item(253, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(253, 16).Value = 2000

# Row 254
$ws.Cells.Item(254, 4).Value = 44488
$ws.Cells.Item(254, 10).Value = 85
$ws.Cells.Item(254, 11).Value = 2000
$ws.Cells.Item(254, 12).Value = 2000
$ws.Cells.Item(254, 13).Value = 2000
$ws.Cells.Item(254, 15).Value = "Región Metropolitana"
$ws.Cells.Item(254, 16).Value = 667

# Row 255
$ws.Cells.Item(255, 4).Value = 44369
$ws.Cells.Item(255, 10).Value = 30

# Row 256
$ws.Cells.Item(256, 4).Value = 44596
$ws.Cells.Item(256, 10).Value = 20
$ws.Cells.Item(256, 11).Value = 5000
$ws.Cells.Item(256, 12).Value = 5000
$ws.Cells.Item(256, 13).Value = 5000
$ws.Cells.Item(256, 16).Value = 1667

# Row 257
$ws.Cells.Item(257, 4).Value = 44399
$ws.Cells.Item(257, 10).Value = 30
$ws.Cells.Item(257, 11).Value = 10000
$ws.Cells.Item(257, 12).Value = 10000
$ws.Cells.Item(257, 13).Value = 10000
$ws.Cells.Item(257, 16).Value = 3333

# Row 258
$ws.Cells.Item(258, 4).Value = 44425
$ws.Cells.Item(258, 11).Value = 8000
$ws.Cells.Item(258, 12).Value = 8000
$ws.Cells.Item(258, 13).Value = 8000
$ws.Cells.Item(258, 16).Value = 2667

# Row 259
$ws.Cells.Item(259, 4).Value = 44512
$ws.Cells.Item(259, 10).Value = 20
$ws.Cells.Item(259, 11).Value = 5000
$ws.Cells.Item(259, 12).Value = 5000
$ws.Cells.Item(259, 13).Value = 5000
$ws.Cells.Item(259, 16).Value = 1667

# Row 260
$ws.Cells.Item(260, 4).Value = 44397
$ws.Cells.Item(260, 10).Value = 10
$ws.Cells.Item(260, 11).Value = 7000
$ws.Cells.Item(260, 12).Value = 7000
$ws.Cells.Item(260, 13).Value = 7000
$ws.Cells.Item(260, 16).Value = 2333

# Row 261
$ws.Cells.Item(261, 4).Value = 44181
$ws.Cells.Item(261, 10).Value = 75
$ws.Cells.Item(261, 11).Value = 5000
$ws.Cells.Item(261, 12).Value = 5000
$ws.Cells.Item(261, 13).Value = 5000
$ws.Cells.Item(261, 16).Value = 1667

# Row 262
$ws.Cells.Item(262, 4).Value = 44747
$ws.Cells.Item(262, 10).Value = 55
$ws.Cells.Item(262, 11).Value = 9000
$ws.Cells.Item(262, 12).Value = 10000
$ws.Cells.Item(262, 13).Value = 9636
$ws.Cells.Item(262, 16).Value = 3212

# Row 263
$ws.Cells.Item(263, 4).Value = 44357
$ws.Cells.Item(263, 11).Value = 5000
$ws.Cells.Item(263, 12).Value = 5000
$ws.Cells.Item(263, 13).Value = 5000
$ws.Cells.Item(263, 16).Value = 1667

# Row 264
$ws.Cells.Item(264, 4).Value = 44757
$ws.Cells.Item(264, 10).Value = 40

# Row 265
$ws.Cells.Item(265, 4).Value = 44321
$ws.Cells.Item(265, 10).Value = 45
$ws.Cells.Item(265, 11).Value = 6000
$ws.Cells.Item(265, 12).Value = 6000
$ws.Cells.Item(265, 13).Value = 6000
$ws.Cells.Item(265, 16).Value = 2000

# Row 266
$ws.Cells.Item(266, 4).Value = 44438
$ws.Cells.Item(266, 11).Value = 10000
$ws.Cells.Item(266, 12).Value = 10000
$ws.Cells.Item(266, 13).Value = 10000
$ws.Cells.Item(266, 16).Value = 3333

# Row 267
$ws.Cells.Item(267, 4).Value = 44355
$ws.Cells.Item(267, 10).Value = 20
$ws.Cells.Item(267, 11).Value = 5000
$ws.Cells.Item(267, 12).Value = 5000
$ws.Cells.Item(267, 13).Value = 5000
$ws.Cells.Item(267, 16).Value = 1667

# Row 268
$ws.Cells.Item(268, 4).Value = 44391
$ws.Cells.Item(268, 10).Value = 55
$ws.Cells.Item(268, 11).Value = 7000
$ws.Cells.Item(268, 12).Value = 7000
$ws.Cells.Item(268, 13).Value = 7000
$ws.Cells.Item(268, 16).Value = 2333

# Row 269
$ws.Cells.Item(269, 4).Value = 44186
$ws.Cells.Item(269, 12).Value = 4000
$ws.Cells.Item(269, 13).Value = 4000
$ws.Cells.Item(269, 16).Value = 1333

# Row 270
$ws.Cells.Item(270, 4).Value = 44189
$ws.Cells.Item(270, 11).Value = 4000
$ws.Cells.Item(270, 13).Value = 4500
$ws.Cells.Item(270, 16).Value = 1500

# Row 271
$ws.Cells.Item(271, 4).Value = 44609
$ws.Cells.Item(271, 10).Value = 40

# Row 272 (new)
$ws.Cells.Item(272, 1).Value = 10
$ws.Cells.Item(272, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(272, 3).Value = "La Araucanía"
$ws.Cells.Item(272, 4).Value = 44358
$ws.Cells.Item(272, 5).Value = 9
$ws.Cells.Item(272, 6).Value = 100112039
$ws.Cells.Item(272, 7).Value = "Ciboulette"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 30
$ws.Cells.Item(272, 11).Value = 5000
$ws.Cells.Item(272, 12).Value = 5000
$ws.Cells.Item(272, 13).Value = 5000
$ws.Cells.Item(272, 14).Value = "`$/docena de atados"
$ws.Cells.Item(272, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(272, 16).Value = 1667
$ws.Cells.Item(272, 17).Value = 3
$ws.Cells.Item(272, 18).Value = "Hortaliza"
$ws.Cells.Item(272, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
